$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.103.01'
$ws.Range("E2").Value = '  +5.68%  '

$ws.Range("D3").Value = '3.121.76'
$ws.Range("E3").Value = '  +3.77%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'" + '584.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.50%  '

$ws.Range("D6").Value = "'" + '145.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.51%  '

$ws.Range("D7").Value = "'" + '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.113.21'
$ws.Range("E8").Value = '  +3.85%  '

$ws.Range("E9").Value = '  +1.69%  '

$ws.Range("E10").Value = '  +12.23%  '

$ws.Range("E11").Value = '  +7.35%  '

$ws.Range("D12").Value = "'" + '0.468'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.32%  '

$ws.Range("D13").Value = "'" + '0.0000250'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.59%  '

$ws.Range("D14").Value = "'" + '35.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.11%  '

$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '3.638.65'
$ws.Range("E16").Value = '  +3.79%  '

$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '63.017.58'
$ws.Range("E18").Value = '  +5.57%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.111.21'
$ws.Range("E19").Value = '  +3.51%  '

$ws.Range("D20").Value = "'" + '468.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.89%  '

$ws.Range("D21").Value = "'" + '14.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.24%  '

$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").Value = "'" + '7.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.10%  '

$ws.Range("D24").Value = "'" + '13.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.59%  '

$ws.Range("D25").Value = "'" + '82.17'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("E28").Value = '  +5.01%  '

$ws.Range("D29").Value = "'" + '8.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.13%  '

$ws.Range("D31").Value = "'" + '6.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.95%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'" + '27.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.37%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'" + '0.111'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.91%  '

$ws.Range("D34").Value = '0.0₃0877'
$ws.Range("E34").Value = '  +10.89%  '

$ws.Range("E35").Value = '  +13.90%  '

$ws.Range("E36").Value = '  +3.78%  '

$ws.Range("E37").Value = '  +2.36%  '

$ws.Range("D38").Value = "'" + '3.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.62%  '

$ws.Range("D39").Value = "'" + '50.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.56%  '

$ws.Range("D40").Value = "'" + '434.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.57%  '

$ws.Range("D41").Value = "'" + '8.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.89%  '

$ws.Range("D42").Value = '2.939.29'
$ws.Range("E42").Value = '  +6.19%  '

$ws.Range("D43").Value = "'" + '0.0370'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.75%  '

$ws.Range("E44").Value = '  +10.21%  '

$ws.Range("E45").Value = '  +4.31%  '

$ws.Range("D46").Value = "'" + '2.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.05%  '

$ws.Range("D47").Value = "'" + '35.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.52%  '

$ws.Range("D49").Value = "'" + '123.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("E50").Value = '  +0.81%  '

$ws.Range("D51").Value = "'" + '24.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.16%  '
